# feat: add 2022-Q1 data
#
# - Insert a new worksheet "2022-Q1" positioned between "2021-Q4" and "总计",
#   populated with the quarter's fund holdings (same column layout as the
#   "2021-Q4" sheet).
# - Insert a new top row into "总计" summarising the new quarter
#   (date / holdings count / holding value), pushing the existing rows down.

function Set-TextCell($range, [string]$value) {
    # Forces the cell to store a STRING even when $value looks numeric
    # (e.g. "000800" or "4.31"), matching Excel's "format as Text, then
    # type the value" behaviour -- without leaving a lingering explicit
    # number-format style on the cell once we're done.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$wb = $excel.ActiveWorkbook

$q4Sheet   = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------
# 1. New "2022-Q1" sheet: clone "2021-Q4" (same header/style layout),
#    rename, trim to the 3 funds reported this quarter, then overwrite
#    the data cells with the new quarter's figures.
# ---------------------------------------------------------------------
$q4Sheet.Copy($totalSheet)
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

# Only 3 funds this quarter (template sheet had 6 rows of data: 2-7).
$newSheet.Rows("5:7").Delete()

# Row 2 - 000800 华商未来主题混合
Set-TextCell $newSheet.Range("B2") "000800"
Set-TextCell $newSheet.Range("C2") "华商未来主题混合"
Set-TextCell $newSheet.Range("D2") "4.31"
Set-TextCell $newSheet.Range("E2") "84.71"
Set-TextCell $newSheet.Range("F2") "4.25"
Set-TextCell $newSheet.Range("G2") "0.1832"
$newSheet.Range("H2").Value = 4

# Row 3 - 001449 华商双驱优选灵活配置混合
Set-TextCell $newSheet.Range("B3") "001449"
Set-TextCell $newSheet.Range("C3") "华商双驱优选灵活配置混合"
Set-TextCell $newSheet.Range("D3") "2.88"
Set-TextCell $newSheet.Range("E3") "92.39"
Set-TextCell $newSheet.Range("F3") "5.00"
Set-TextCell $newSheet.Range("G3") "0.1440"
$newSheet.Range("H3").Value = 2

# Row 4 - 010403 华商景气优选混合
Set-TextCell $newSheet.Range("B4") "010403"
Set-TextCell $newSheet.Range("C4") "华商景气优选混合"
Set-TextCell $newSheet.Range("D4") "0.61"
Set-TextCell $newSheet.Range("E4") "84.95"
Set-TextCell $newSheet.Range("F4") "4.17"
Set-TextCell $newSheet.Range("G4") "0.0254"
$newSheet.Range("H4").Value = 8

# ---------------------------------------------------------------------
# 2. "总计" sheet: insert a new row right under the header for 2022-Q1,
#    pushing the previous 2021-Q4 / 2021-Q3 rows down by one.
#    (Re-fetch by name: inserting/copying sheets shifts tab positions, so
#    an index-based handle captured earlier may now point elsewhere.)
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows("2:2").Insert()
$totalSheet.Range("A2:D2").ClearFormats()

$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)   # xlPasteFormats - reuse the existing style

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 3
$totalSheet.Range("D2").Value = 0.35

# The "A" column is a 0-based row index, not the shifted rows' original
# values -- renumber the two pushed-down rows to match.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
